# Recollected Jfreechart and Jodatime results and updated Results.xlsx
#
# On the "Selection" sheet, the manual-tests table (rows 11-12, the
# JFreeChart / Jodatime rows under "Test Selection Manual Tests") had no
# recorded values yet. Fill in the recollected results: all zero except
# Jodatime's "Relative Function" count, which is 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Selection")

# Row 11 = JFreeChart, Row 12 = Jodatime; columns B..I are the 8 metrics.
$row11 = @(0, 0, 0, 0, 0, 0, 0, 0)
$row12 = @(0, 0, 0, 0, 0, 0, 1, 0)

$cols = @("B", "C", "D", "E", "F", "G", "H", "I")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $row11[$i]
    $ws.Range($cols[$i] + "12").Value = $row12[$i]
}

$ws.Range("G18").Select()
